$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 21 ("День самостоятельной работы") with the pair number and
# placeholder dashes for room / position / teacher / pair type.
$ws.Range("C21").Value = 1
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "-"
$ws.Range("G21").Value = "-"
$ws.Range("H21").Value = "-"

# Leave the selection where the author left it after editing row 21.
[void]$ws.Range("H22").Select()
